$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5036.593
$ws.Range("I62").Value = 3781.0454
$ws.Range("J62").Value = 10561
$ws.Range("K62").Value = 3781.0454
$ws.Range("L62").Value = 10561
$ws.Range("M62").Value = -3157.0454
$ws.Range("N62").Value = -11809
$ws.Range("H65").Value = 5036.593
$ws.Range("I65").Value = 3781.0454
$ws.Range("J65").Value = 10561
$ws.Range("K65").Value = 18905.227
$ws.Range("L65").Value = 52805
$ws.Range("M65").Value = -15785.227
$ws.Range("N65").Value = -59045
$ws.Range("H107").Value = 744.25
$ws.Range("I107").Value = 775.8
$ws.Range("J107").Value = 691.6667
$ws.Range("K107").Value = 775.8
$ws.Range("L107").Value = 691.6667
$ws.Range("M107").Value = 1144.2
$ws.Range("N107").Value = -4531.6667
$ws.Range("H111").Value = 1492.6666
$ws.Range("J111").Value = 1500
$ws.Range("L111").Value = 4500
$ws.Range("N111").Value = -10634
$ws.Range("H112").Value = 1379.9375
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H132").Value = 2192.1304
$ws.Range("I132").Value = 2258.1904
$ws.Range("J132").Value = 1498.5
$ws.Range("K132").Value = 6774.5712
$ws.Range("L132").Value = 4495.5
$ws.Range("M132").Value = -4244.5712
$ws.Range("N132").Value = -9555.5
$ws.Range("H137").Value = 213165.1
$ws.Range("I137").Value = 2707.1785
$ws.Range("J137").Value = 1119753.1
$ws.Range("K137").Value = 8121.5355
$ws.Range("L137").Value = 3359259.3
$ws.Range("M137").Value = -5571.5355
$ws.Range("N137").Value = -3364359.3
$ws.Range("H138").Value = 2211.1667
$ws.Range("J138").Value = 3571
$ws.Range("L138").Value = 10713
$ws.Range("N138").Value = -20993

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1228.3
$ws.Range("I2").Value = 1252.2222
$ws.Range("K2").Value = 1252.2222
$ws.Range("M2").Value = -1139.2222
$ws.Range("H32").Value = 7109.9043
$ws.Range("J32").Value = 26272
$ws.Range("L32").Value = 26272
$ws.Range("N32").Value = -26846
$ws.Range("H45").Value = 6952343
$ws.Range("I45").Value = 9512.5
$ws.Range("J45").Value = 31252250
$ws.Range("K45").Value = 9512.5
$ws.Range("L45").Value = 31252250
$ws.Range("M45").Value = -9135.5
$ws.Range("N45").Value = -31253004
$ws.Range("H63").Value = 4724.75
$ws.Range("I63").Value = 2299.6667
$ws.Range("J63").Value = 12000
$ws.Range("K63").Value = 2299.6667
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -1613.6667
$ws.Range("N63").Value = -13372
$ws.Range("H66").Value = 4724.75
$ws.Range("I66").Value = 2299.6667
$ws.Range("J66").Value = 12000
$ws.Range("K66").Value = 11498.3335
$ws.Range("L66").Value = 60000
$ws.Range("M66").Value = -8066.333500000001
$ws.Range("N66").Value = -66864
$ws.Range("H74").Value = 2072.3823
$ws.Range("I74").Value = 1421
$ws.Range("J74").Value = 3002.9285
$ws.Range("K74").Value = 1421
$ws.Range("L74").Value = 3002.9285
$ws.Range("M74").Value = -547
$ws.Range("N74").Value = -4750.9285
$ws.Range("H77").Value = 2072.3823
$ws.Range("I77").Value = 1421
$ws.Range("J77").Value = 3002.9285
$ws.Range("K77").Value = 7105
$ws.Range("L77").Value = 15014.6425
$ws.Range("M77").Value = -2737
$ws.Range("N77").Value = -23750.6425
$ws.Range("H116").Value = 1228.3
$ws.Range("I116").Value = 1252.2222
$ws.Range("K116").Value = 1252.2222
$ws.Range("M116").Value = 1041.7778
$ws.Range("H122").Value = 2178.6072
$ws.Range("I122").Value = 2151.524
$ws.Range("K122").Value = 6454.572
$ws.Range("M122").Value = -4004.572
$ws.Range("H132").Value = 2398.3333
$ws.Range("I132").Value = 1526.7222
$ws.Range("J132").Value = 4141.5557
$ws.Range("K132").Value = 4580.1666
$ws.Range("L132").Value = 12424.6671
$ws.Range("M132").Value = -2050.1666
$ws.Range("N132").Value = -17484.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1228.3
$ws.Range("I3").Value = 1252.2222
$ws.Range("K3").Value = 1252.2222
$ws.Range("M3").Value = -1138.2222
$ws.Range("H20").Value = 107155
$ws.Range("I20").Value = 144221.5
$ws.Range("J20").Value = 3368.8
$ws.Range("K20").Value = 144221.5
$ws.Range("L20").Value = 3368.8
$ws.Range("M20").Value = -143974.5
$ws.Range("N20").Value = -3862.8
$ws.Range("H105").Value = 41536.52
$ws.Range("I105").Value = 59625.59
$ws.Range("K105").Value = 59625.59
$ws.Range("M105").Value = -57878.59

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 635.8570999999999
$ws.Range("I22").Value = 658.5
$ws.Range("K22").Value = 658.5
$ws.Range("M22").Value = -308.5
$ws.Range("H58").Value = 1627.6052
$ws.Range("I58").Value = 1195.8518
$ws.Range("J58").Value = 2687.3635
$ws.Range("K58").Value = 1195.8518
$ws.Range("L58").Value = 2687.3635
$ws.Range("M58").Value = -992.8517999999999
$ws.Range("N58").Value = -3093.3635
$ws.Range("H122").Value = 2239.72
$ws.Range("I122").Value = 1783
$ws.Range("J122").Value = 4637.5
$ws.Range("K122").Value = 5349
$ws.Range("L122").Value = 13912.5
$ws.Range("M122").Value = -2899
$ws.Range("N122").Value = -18812.5
$ws.Range("H132").Value = 1969.6
$ws.Range("I132").Value = 1514.75
$ws.Range("K132").Value = 4544.25
$ws.Range("M132").Value = -2014.25
$ws.Range("H134").Value = 5942.1816
$ws.Range("I134").Value = 6540.4443
$ws.Range("J134").Value = 3250
$ws.Range("K134").Value = 19621.3329
$ws.Range("L134").Value = 9750
$ws.Range("M134").Value = -17086.3329
$ws.Range("N134").Value = -14820
$ws.Range("H136").Value = 1627.6052
$ws.Range("I136").Value = 1195.8518
$ws.Range("J136").Value = 2687.3635
$ws.Range("K136").Value = 3587.5554
$ws.Range("L136").Value = 8062.0905
$ws.Range("M136").Value = -1037.5554
$ws.Range("N136").Value = -13162.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 216
$ws.Range("I50").Value = 102.4
$ws.Range("K50").Value = 307.2
$ws.Range("M50").Value = 173.8
$ws.Range("H53").Value = 216
$ws.Range("I53").Value = 102.4
$ws.Range("K53").Value = 307.2
$ws.Range("M53").Value = 173.8
$ws.Range("H131").Value = 1808
$ws.Range("I131").Value = 1026.2858
$ws.Range("J131").Value = 2589.7144
$ws.Range("K131").Value = 3078.8574
$ws.Range("L131").Value = 7769.1432
$ws.Range("M131").Value = 1961.1426
$ws.Range("N131").Value = -17849.1432
$ws.Range("H140").Value = 2918.6875
$ws.Range("I140").Value = 2333.1667
$ws.Range("K140").Value = 6999.500100000001
$ws.Range("M140").Value = -1819.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H110").Value = 92381.55
$ws.Range("J110").Value = 92381.55
$ws.Range("L110").Value = 92381.55
$ws.Range("N110").Value = -100561.55
$ws.Range("H126").Value = 2674.625
$ws.Range("I126").Value = 2281.2173
$ws.Range("K126").Value = 6843.651899999999
$ws.Range("M126").Value = -4373.651899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 18066.166
$ws.Range("I46").Value = 34133.668
$ws.Range("J46").Value = 1998.6666
$ws.Range("K46").Value = 34133.668
$ws.Range("L46").Value = 1998.6666
$ws.Range("M46").Value = -33945.668
$ws.Range("N46").Value = -2374.6666
$ws.Range("H61").Value = 3600
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H96").Value = 88566.664
$ws.Range("J96").Value = 88566.664
$ws.Range("L96").Value = 88566.664
$ws.Range("N96").Value = -94058.664
$ws.Range("H112").Value = 34177.4
$ws.Range("J112").Value = 34177.4
$ws.Range("L112").Value = 34177.4
$ws.Range("N112").Value = -37131.4
$ws.Range("H113").Value = 3600
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 12521501
$ws.Range("I122").Value = 27335.666
$ws.Range("K122").Value = 82006.99800000001
$ws.Range("M122").Value = -79556.99800000001
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 13557.846
$ws.Range("I132").Value = 24726.5
$ws.Range("K132").Value = 74179.5
$ws.Range("M132").Value = -71649.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 49999
$ws.Range("I21").Value = 49999
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 49999
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -49764
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 49999
$ws.Range("I35").Value = 49999
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 49999
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -49709
$ws.Range("N35").ClearContents()
$ws.Range("H70").Value = 250000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 250000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 250000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -250630
$ws.Range("H73").Value = 250000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 250000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 250000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -252184
$ws.Range("H107").Value = 3473.2
$ws.Range("I107").Value = 2900.3333
$ws.Range("J107").Value = 4332.5
$ws.Range("K107").Value = 8700.999899999999
$ws.Range("L107").Value = 12997.5
$ws.Range("M107").Value = -6780.999899999999
$ws.Range("N107").Value = -16837.5
$ws.Range("H126").Value = 2011.9143
$ws.Range("I126").Value = 1680.6086
$ws.Range("K126").Value = 5041.825800000001
$ws.Range("M126").Value = -2571.825800000001
$ws.Range("H130").Value = 67000
$ws.Range("J130").Value = 67000
$ws.Range("L130").Value = 67000
$ws.Range("N130").Value = -77040
